$d = $word.ActiveDocument

# The document contains two tables; the second one is
# "Analisi minacce e controlli" (Minaccia / Probabilità / Controllo / Fattibilità).
$table = $d.Tables.Item(2)

function Set-CellText {
    param([int]$Row, [int]$Col, [string]$Text, [bool]$Center)
    $cell = $table.Cell($Row, $Col)
    $range = $cell.Range
    # Trim trailing paragraph/cell markers from the range so we only
    # replace the actual text content of the cell.
    $range.End = $range.End - 1
    $range.Text = $Text
    if ($Center) {
        $range.ParagraphFormat.Alignment = 1
    }
}

# Row 2: Furto di identità(utente)
Set-CellText 2 2 "Alta" $true
Set-CellText 2 3 "Log delle operazioni con etichetta utente." $false
Set-CellText 2 4 "Basso costo e trasparente." $false

# Row 3: Compromissione codice(bagaglio)
Set-CellText 3 2 "Media" $true
Set-CellText 3 3 "Assegnazione codice secondo una precisa struttura alfanumerica." $false
Set-CellText 3 4 "Medio costo, evita la clonazione del codice bagaglio." $false

# Row 4: Alterazione posizione(bagaglio)
Set-CellText 4 2 "Media" $true
Set-CellText 4 3 "Utilizzo del sensore in condizioni ottimali che non compromettano la corretta recezione del segnale." $false
Set-CellText 4 4 "Alto costo, vincolare l’utilizzo del sensore in situazioni che escludano l’alta quota." $false

# Row 5: previously fully empty -> now "Crash del sito web"
Set-CellText 5 1 "Crash del sito web" $false
Set-CellText 5 2 "Bassa" $true
Set-CellText 5 3 "Creazione di un sito che sia in grado di gestire n richieste contemporaneamente." $false
Set-CellText 5 4 "Alto costo, ottimizzare database ed effettuare backup frequente." $false
